$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted. It belongs right above the
# current row 272, so push that row (and everything below it) down by one
# row, then fill the freshly inserted row with the new record's data.
$ws.Rows.Item(272).Insert()

$ws.Cells.Item(272, 1).Value = 9
$ws.Cells.Item(272, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(272, 3).Value = "Metropolitana"
$ws.Cells.Item(272, 4).Value = 44578
$ws.Cells.Item(272, 5).Value = 13
$ws.Cells.Item(272, 6).Value = 100112032
$ws.Cells.Item(272, 7).Value = "Zapallo italiano"
$ws.Cells.Item(272, 8).Value = "Sin especificar"
$ws.Cells.Item(272, 9).Value = "Primera"
$ws.Cells.Item(272, 10).Value = 61
$ws.Cells.Item(272, 11).Value = 13000
$ws.Cells.Item(272, 12).Value = 14000
$ws.Cells.Item(272, 13).Value = 13508
$ws.Cells.Item(272, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(272, 15).Value = "Región del Maule"
$ws.Cells.Item(272, 16).Value = 225
$ws.Cells.Item(272, 17).Value = 60
$ws.Cells.Item(272, 18).Value = "Hortaliza"
